$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - reorder columns B/C, G/H, I/J
$ws.Range("B1").Value = "expire_date"
$ws.Range("C1").Value = "recent_location"
$ws.Range("G1").Value = "origin_location"
$ws.Range("H1").Value = "gender"
$ws.Range("I1").Value = "issue_place"
$ws.Range("J1").Value = "issue_date"

# Data row (row 2) - new record values
# Leading-apostrophe forces text (preserves leading zero); reset Style
# afterward so the cell keeps the same (default/no) style as the source.
$ws.Range("A2").Value = "'049300009355"
$ws.Range("A2").Style = "Normal"
$ws.Range("B2").Value = "14/05/2025"
$ws.Range("C2").Value = "Trung Toàn 1 Tam Quang, Núi Thành, Quảng Nam"
$ws.Range("D2").Value = "ĐINH THỊ TÚ TRANG"
$ws.Range("E2").Value = "14/05/2000"
$ws.Range("F2").Value = "Việt Nam"
$ws.Range("G2").Value = "Tam Quang, Núi Thành, Quảng Nam"
$ws.Range("H2").Value = "Nữ"
$ws.Range("I2").Value = "CỤC TRƯỞNG CỤC CẢNH SÁT QUẢN LÝ HÀNH CHÍNH VỀ TRẬT TỰ XÃ HỘI"
$ws.Range("J2").Value = "15/08/2021"
